$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin/Link/Price/Volume text updates (values that Excel will not mis-parse as numbers) ---
$textUpdates = @{
    "D2" = '96.117.03'
    "E2" = '  +0.38%  '
    "D3" = '3.544.83'
    "E3" = '  -1.81%  '
    "E4" = '  +0.12%  '
    "E5" = '  -0.01%  '
    "E6" = '  -0.55%  '
    "E7" = '  +10.57%  '
    "E8" = '  +0.00%  '
    "E9" = '  +6.27%  '
    "E10" = '  +0.15%  '
    "D11" = '3.546.52'
    "E11" = '  -1.74%  '
    "E12" = '  +0.58%  '
    "E13" = '  +0.56%  '
    "E14" = '  -0.07%  '
    "D15" = '4.213.57'
    "E15" = '  -1.50%  '
    "D16" = '95.996.73'
    "E16" = '  +0.30%  '
    "E17" = '  +1.75%  '
    "D18" = '3.547.15'
    "E18" = '  -1.66%  '
    "E19" = '  -2.72%  '
    "E20" = '  -2.69%  '
    "E21" = '  -2.38%  '
    "E22" = '  +8.57%  '
    "E23" = '  -1.18%  '
    "E24" = '  -6.72%  '
    "E25" = '  +2.24%  '
    "E26" = '  +0.60%  '
    "E27" = '  -1.40%  '
    "E28" = '  -0.85%  '
    "D29" = '3.742.88'
    "E29" = '  -1.54%  '
    "E30" = '  +6.98%  '
    "E31" = '  -5.27%  '
    "E32" = '  -0.88%  '
    "E33" = '  +0.01%  '
    "E34" = '  +1.93%  '
    "E35" = '  -0.01%  '
    "E36" = '  -2.76%  '
    "B37" = 'RenderToken'
    "C37" = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    "E37" = '  +5.24%  '
    "B38" = 'Bittensor'
    "C38" = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    "E38" = '  +5.66%  '
    "E39" = '  -0.32%  '
    "E40" = '  +7.03%  '
    "E41" = '  +0.00%  '
    "E42" = '  -0.87%  '
    "E43" = '  -3.38%  '
    "E44" = '  +4.06%  '
    "E45" = '  -1.13%  '
    "B46" = 'VeChain'
    "C46" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "E46" = '  +0.69%  '
    "B47" = 'Stacks'
    "C47" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    "E47" = '  +0.15%  '
    "B48" = 'WhiteBITCoin'
    "C48" = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    "E48" = '  -1.09%  '
    "E49" = '  -3.50%  '
    "E50" = '  -0.15%  '
    "E51" = '  +0.04%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Price cells whose new values look like plain numbers (e.g. "1.64", "0.0420") ---
# Force the cell to Text format first so Excel stores the exact original string
# (keeping trailing zeros / leading zeros intact) instead of silently converting it
# to a numeric value, then restore the default (unstyled) cell style afterwards.
$numericLookingUpdates = @{
    "D5" = '239.48'
    "D6" = '651.91'
    "D7" = '1.64'
    "D8" = '0.406'
    "D9" = '1.07'
    "D12" = '43.32'
    "D14" = '6.34'
    "D17" = '0.0000260'
    "D19" = '8.02'
    "D20" = '12.35'
    "D21" = '17.61'
    "D22" = '0.535'
    "D23" = '504.84'
    "D24" = '3.38'
    "D25" = '6.82'
    "D27" = '95.63'
    "D28" = '12.65'
    "D30" = '0.150'
    "D32" = '11.31'
    "D35" = '1.00'
    "D36" = '31.14'
    "D37" = '8.67'
    "D38" = '607.23'
    "D40" = '1.59'
    "D42" = '0.149'
    "D43" = '0.894'
    "D44" = '1.80'
    "D45" = '5.67'
    "D46" = '0.0420'
    "D47" = '2.26'
    "D48" = '23.51'
    "D49" = '33.06'
    "D51" = '8.12'
}
foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
    $cell.Style = "Normal"
}
